$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: clear out the two old long Cypher query strings that currently
# live in A2 and B2. Doing this first (before inserting the new column)
# frees their shared-string slots and lets them be reclaimed/compacted, so
# that after we write the new content the shared-strings table ends up in
# the same order as the target workbook (existing strings keep their slot,
# freed ones are dropped, new ones are appended at the end).
$ws.Range("A2").Value2 = "placeholder1"
$ws.Range("B2").Value2 = "placeholder2"

# --- Step 2: insert a brand new column A ("TabName" / "CasesTab"), shifting
# the previous columns A-D to B-E. The existing column widths (which already
# match the target widths for the shifted-right columns) move along with
# the columns automatically.
$ws.Columns("A").Insert()

# --- Step 3: populate the new column A
$ws.Range("A1").Value2 = "TabName"
$ws.Range("A2").Value2 = "CasesTab"

# --- Step 4: write the new query text into B2/C2 (which is where the old
# A2/B2 placeholders ended up after the column insert). These replace the
# placeholders we set in Step 1.
$casesQuery = @"
MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.gender = "MALE"
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '') AS ``Case ID``,
    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,
    COALESCE(a.arm_id, '') AS ``Arm``,
    COALESCE(a.arm_drug, '') AS ``Arm Treatment``,
    COALESCE(c.disease, '') AS ``Diagnosis``,
    COALESCE(c.gender, '') AS ``Gender``,
    COALESCE(c.race, '') AS ``Race``,
    COALESCE(c.ethnicity, '') AS ``Ethnicity``
"@

$statQuery = @"
MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.gender = "MALE"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials
"@

$ws.Range("B2").Value2 = $casesQuery
$ws.Range("C2").Value2 = $statQuery

# --- Step 5: formatting to match the updated layout
# New column A is a narrow "best fit"-style column for the short tab name.
$ws.Columns("A").ColumnWidth = 8.0

# Row 2 now needs to be tall enough to show the much longer, multi-line
# query text.
$ws.Rows(2).RowHeight = 174

# Update the active selection to reflect where the author left off editing.
$ws.Range("C4").Select()
